# This script reproduces the commit "Fruta / hortaliza, semanal":
# two new weekly price records are inserted at the top of the
# "Pepino ensalada" data block (rows 297-298), pushing the existing
# records down by two rows (old 297..360 -> new 299..362).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 297 (existing rows shift down, carrying
# their formatting/date style with them).
$ws.Rows.Item(297).Insert()
$ws.Rows.Item(298).Insert()

# --- New row 297 ---
$ws.Range("A297").Value = 5
$ws.Range("B297").Value = "Macroferia Regional de Talca"
$ws.Range("C297").Value = "Maule"
$ws.Range("D297").Value = 44637
$ws.Range("E297").Value = 7
$ws.Range("F297").Value = 100112043
$ws.Range("G297").Value = "Pepino ensalada"
$ws.Range("H297").Value = "Sin especificar"
$ws.Range("I297").Value = "Primera"
$ws.Range("J297").Value = 300
$ws.Range("K297").Value = 18000
$ws.Range("L297").Value = 18000
$ws.Range("M297").Value = 18000
$ws.Range("N297").Value = "$/caja 60 unidades"
$ws.Range("O297").Value = "Región de Arica y Parinacota"
$ws.Range("P297").Value = 300
$ws.Range("Q297").Value = 60
$ws.Range("R297").Value = "Hortaliza"

# --- New row 298 ---
$ws.Range("A298").Value = 5
$ws.Range("B298").Value = "Macroferia Regional de Talca"
$ws.Range("C298").Value = "Maule"
$ws.Range("D298").Value = 44637
$ws.Range("E298").Value = 7
$ws.Range("F298").Value = 100112043
$ws.Range("G298").Value = "Pepino ensalada"
$ws.Range("H298").Value = "Sin especificar"
$ws.Range("I298").Value = "Primera"
$ws.Range("J298").Value = 300
$ws.Range("K298").Value = 20000
$ws.Range("L298").Value = 20000
$ws.Range("M298").Value = 20000
$ws.Range("N298").Value = "$/caja 80 unidades"
$ws.Range("O298").Value = "Región del Maule"
$ws.Range("P298").Value = 250
$ws.Range("Q298").Value = 80
$ws.Range("R298").Value = "Hortaliza"
